$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update Solar (column E) for 2022 (row 24): 586 -> 587
$ws.Range("E24").Value = 587

# Update Energy Storage (column C) for 2024 (row 26): 3 -> 8
$ws.Range("C26").Value = 8

# Update Solar (column E) for 2024 (row 26): 678 -> 949
$ws.Range("E26").Value = 949
